$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 96 (shifts existing rows 96:211 down to 97:212)
$ws.Rows(96).Insert()

# Populate the newly inserted row with the new price-report record
$ws.Cells.Item(96,1).Value  = 3
$ws.Cells.Item(96,2).Value  = "Femacal de La Calera"
$ws.Cells.Item(96,3).Value  = "Coquimbo"
$ws.Cells.Item(96,4).Value  = 44483
$ws.Cells.Item(96,5).Value  = 5
$ws.Cells.Item(96,6).Value  = 100112043
$ws.Cells.Item(96,7).Value  = "Pepino ensalada"
$ws.Cells.Item(96,8).Value  = "Sin especificar"
$ws.Cells.Item(96,9).Value  = "Primera"
$ws.Cells.Item(96,10).Value = 105
$ws.Cells.Item(96,11).Value = 14000
$ws.Cells.Item(96,12).Value = 14500
$ws.Cells.Item(96,13).Value = 14262
$ws.Cells.Item(96,14).Value = "$/caja 70 unidades"
$ws.Cells.Item(96,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(96,16).Value = 204
$ws.Cells.Item(96,17).Value = 70
$ws.Cells.Item(96,18).Value = "Hortaliza"
